$wb = $excel.ActiveWorkbook

# --- Sheet "Uitvallers": append three more riders who abandoned the race ---
$wsUitvallers = $wb.Worksheets.Item("Uitvallers")
$wsUitvallers.Activate()

$wsUitvallers.Range("A22").Value = "Carlos Rodriguez"
$wsUitvallers.Range("A23").Value = "Cyril Barthe"
$wsUitvallers.Range("A24").Value = "Enric Mas"

$wsUitvallers.Range("A25").Select()

# --- Sheet "Huidig": fill in column S (stage-18 predictions) for rows 6-19 ---
$wsHuidig = $wb.Worksheets.Item("Huidig")
$wsHuidig.Activate()

$wsHuidig.Range("S6").Value  = "Ben O'Connor"
$wsHuidig.Range("S7").Value  = "Tadej Pogacar"
$wsHuidig.Range("S8").Value  = "Jonas Vingegaard"
$wsHuidig.Range("S9").Value  = "Oscar Onley"
$wsHuidig.Range("S10").Value = "Einer Rubio Reyes"
$wsHuidig.Range("S11").Value = "Félix Gall"
$wsHuidig.Range("S12").Value = "Primoz Roglic"
$wsHuidig.Range("S13").Value = "Adam Yates"
$wsHuidig.Range("S14").Value = "Tobias Johannessen"
$wsHuidig.Range("S15").Value = "Sepp Kuss"
$wsHuidig.Range("S16").Value = "Tadej Pogacar"
$wsHuidig.Range("S17").Value = "Jonathan Milan"
$wsHuidig.Range("S18").Value = "Tadej Pogacar"
$wsHuidig.Range("S19").Value = "Florian Lipowitz"

# Update the view of the "Huidig" sheet: scroll/selection moved one column to the right
$excel.ActiveWindow.ScrollColumn = 16
$wsHuidig.Range("S20").Select()

# "Huidig" ends up as the selected/active sheet (instead of "Uitvallers")
